# Generate Report for Handoff
# Renames the handed-off document id from
#   d884ab0e-fede-4adc-95db-8b81b5e70ee7
# to
#   af782309-5a61-44cf-acfc-13ae29cb091c
# across the Overview / zh-cn / de-de sheets, bumps the recorded xliff
# checksum used in the handoff/handback file names, and refreshes the
# handoff/handback timestamps to reflect the new run.

$wb = $excel.ActiveWorkbook

$oldGuid = "d884ab0e-fede-4adc-95db-8b81b5e70ee7"
$newGuid = "af782309-5a61-44cf-acfc-13ae29cb091c"

$oldHash = "a0c62e0ef636e63e59d409f1192de960a31434e2"
$newHash = "684a6c4b50f9301899e07b0a7b40f9f86aea1a87"

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8a71a6dc029c6dccb91a939af8c66cbb599e0465/e2e/$oldGuid.md"

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value2 = "$newGuid.md"

# B2 carries the hyperlink; rebuild it so the `display` text is refreshed
# while the link keeps pointing at the same external target.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $hyperlinkTarget, "", "", "e2e\$newGuid.md")

$wsOverview.Range("G2").Value2 = "2016-09-02 15:12:41"

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $hyperlinkTarget, "", "", "$newGuid.md")

$wsZhCn.Range("G2").Value2 = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value2 = "2016-09-02 15:12:37"

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $hyperlinkTarget, "", "", "$newGuid.md")

$wsDeDe.Range("G2").Value2 = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value2 = "2016-09-02 15:12:41"
